# Updated okcupid.py with function to binarize labels
# This script reproduces the diff: inserts a new block of results
# (essay0-essay0 word/bigram ngram F1 scores) above the existing
# "essay 4" n-gram block, shifts the "essay 4" block + its chart image
# down by 5 rows, changes the second "essay 4" row from Unigram to
# Bigram, and appends a new "essay 4 / sex / trigram" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert 5 new rows above the old row 100 (old rows 100-101,
#    which hold the "essay 4 / Sex / Unigram" rows, shift down to
#    105-106).
# ---------------------------------------------------------------
$ws.Rows("100:104").Insert()

# ---------------------------------------------------------------
# 2. The floating chart picture ("Picture 3") anchored below this
#    block does not automatically follow the inserted rows in this
#    engine, so reposition it to keep the same offsets relative to
#    its anchor rows, which have now shifted from rows 101/109 to
#    rows 106/114 (the picture's xdr anchor is 0-indexed, so it was
#    anchored to the top of row 101 / bottom of row 109, and must
#    now be anchored to the top of row 106 / bottom of row 114).
# ---------------------------------------------------------------
$chart = $ws.Shapes.Item("Picture 3")
$fromOffsetPt = 61451 / 12700.0
$toOffsetPt = 181311 / 12700.0
$newTop = $ws.Rows(106).Top + $fromOffsetPt
$newBottom = $ws.Rows(114).Top + $toOffsetPt
$chart.Top = $newTop
$chart.Height = $newBottom - $newTop

# ---------------------------------------------------------------
# 3. Populate the new rows 100-103 with the "essay0 essay0" n-gram
#    F1-score block, modeled after the existing "Essay 0 all" block
#    in rows 95-98 (copy per-cell formatting, then set the actual
#    values, so columns that stay empty are not touched at all).
# ---------------------------------------------------------------
$ws.Range("B95").Copy($ws.Range("B100"))
$ws.Range("C95").Copy($ws.Range("C100"))
$ws.Range("D95").Copy($ws.Range("D100"))
$ws.Range("E95").Copy($ws.Range("E100"))
$ws.Range("F95").Copy($ws.Range("F100"))
$ws.Range("E96").Copy($ws.Range("E101"))
$ws.Range("B95").Copy($ws.Range("B102"))
$ws.Range("D95").Copy($ws.Range("D102"))
$ws.Range("E95").Copy($ws.Range("E102"))
$ws.Range("F95").Copy($ws.Range("F102"))
$ws.Range("E96").Copy($ws.Range("E103"))

$ws.Range("B100").Value = "essay0 essay0 words"
$ws.Range("C100").Value = "sex"
$ws.Range("D100").Value = "unigram"
$ws.Range("E100").Value = "m F1-Score: 0.6609642301710731"
$ws.Range("F100").Value = 62.1
$ws.Range("B102").Value = "essay0 essay0 bigrams"
$ws.Range("D102").Value = "bigram"
$ws.Range("F102").Value = 58.8

# Shared-string creation order matters for exact reproduction: the
# "f F1-Score" text was authored before the "m F1-Score" text.
$ws.Range("E103").Value = "f F1-Score: 0.5337704918032787"
$ws.Range("E102").Value = "m F1-Score: 0.6310326933056565"

# ---------------------------------------------------------------
# 4. Fix up the shifted "essay 4" rows (now 105-106): row 106's
#    n-gram column changes from Unigram to Bigram.
# ---------------------------------------------------------------
$ws.Range("D106").Value = "Bigram"

# ---------------------------------------------------------------
# 5. Append the new "essay 4 / sex / trigram" row 107 (column E
#    stays untouched/empty, like column G in row 97's pattern).
# ---------------------------------------------------------------
$ws.Range("B106").Copy($ws.Range("B107"))
$ws.Range("C106").Copy($ws.Range("C107"))
$ws.Range("D106").Copy($ws.Range("D107"))
$ws.Range("F106").Copy($ws.Range("F107"))

# Shared-string creation order: "trigram" was authored before
# "essay 4" in the original edit.
$ws.Range("D107").Value = "trigram"
$ws.Range("B107").Value = "essay 4"
$ws.Range("C107").Value = "sex"
$ws.Range("F107").Value = 57.45

# ---------------------------------------------------------------
# 6. Update the view: scroll position and active selection.
# ---------------------------------------------------------------
$ws.Range("I115").Select()
$excel.ActiveWindow.ScrollRow = 85
